# modified Self Evaluation Protocol
# Updates the self-evaluation score cells and refreshes the sheet's
# on-screen selection / scroll position to match the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the evaluation scores -----------------------------------
# "Numbers of Commits in GitHub" score: 88 -> 89
$ws.Range("C9").Value = 89

# "Web Design" score: 15 -> 16
$ws.Range("C11").Value = 16

# C32 holds =SUM(C6:C31) and will recalculate automatically (326 -> 328)

# --- Update the saved view/selection ---------------------------------
# Previously the view was scrolled to A7 with C28 selected; now the
# sheet should show A19 at the top-left with C15 selected.
$ws.Range("C15").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
